# Insert a new worksheet "A06 vie saint gregoire" right after
# "A05 vie sainte dieudonnee" (the new `A` entry for the "vie saint gregoire"
# quote set), which pushes every following sheet's "A##" label up by one.

$wb = $excel.ActiveWorkbook

# Remember whichever sheet is currently active so we can restore the
# selection afterwards -- Worksheets.Add() below activates the new sheet.
$originallyActiveSheet = $wb.ActiveSheet.Name

# --- 1. Rename the trailing sheets (old name -> new name), working backwards
#        from the last sheet so a target name never collides with a
#        not-yet-renamed sheet. ---
$renames = @(
    @("A25 quatre sereurs", "A26 quatre sereurs"),
    @("A23 roy avoit amie", "A24 roy avoit amie"),
    @("A22 jeu des dez", "A23 jeu des dez"),
    @("A21 mauvais riche homme", "A22 mauvais riche homme"),
    @("A20 vieillards tués", "A21 vieillards tués"),
    @("A19 elegy troyes", "A20 elegy troyes"),
    @("A18 richart sans peour", "A19 richart sans peour"),
    @("A17 robert deable", "A18 robert deable"),
    @("A16 guillaume angleterre", "A17 guillaume angleterre"),
    @("A13 vie seint thibault", "A14 vie seint thibault"),
    @("A12 miracle saint servais", "A13 miracle saint servais"),
    @("A11 vie saint sebastien", "A12 vie saint sebastien"),
    @("A10 poines enfer", "A11 poines enfer"),
    @("A09 vie saint leu", "A10 vie saint leu"),
    @("A08 vie glorieux confesseur", "A09 vie glorieux confesseur"),
    @("A07 vie saint jean paulus", "A08 vie saint jean paulus"),
    @("A06 saint jean evangeliste", "A07 saint jean evangeliste")
)

foreach ($pair in $renames) {
    $oldName = $pair[0]
    $newName = $pair[1]
    $sheet = $wb.Worksheets.Item($oldName)
    $sheet.Name = $newName
}

# --- 2. Insert the brand-new sheet right after "A05 vie sainte dieudonnee" ---
$afterSheet = $wb.Worksheets.Item("A05 vie sainte dieudonnee")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "A06 vie saint gregoire"

# --- 3. Populate the new sheet's header row, matching the other sheets. ---
$newSheet.Range("A1").Value = "line_n"
$newSheet.Range("B1").Value = "prev_line"
$newSheet.Range("C1").Value = "line"
$newSheet.Range("D1").Value = "next_line"

# --- 4. Restore the original active-sheet selection. ---
$wb.Worksheets.Item($originallyActiveSheet).Activate()
